$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per commit diff (symbol list refresh).
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "307.02"
Set-TextValue $ws.Range("E2") "0.43%"
Set-TextValue $ws.Range("D3") "39.50"
Set-TextValue $ws.Range("E3") "9.43%"
Set-TextValue $ws.Range("D4") "5.099"
Set-TextValue $ws.Range("E4") "0.78%"
Set-TextValue $ws.Range("D5") "0.08050"
Set-TextValue $ws.Range("E5") "-0.02%"
Set-TextValue $ws.Range("D6") "1.922"
Set-TextValue $ws.Range("E6") "2.74%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "4.193"
Set-TextValue $ws.Range("E7") "1.76%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D8") "7.953"
Set-TextValue $ws.Range("E8") "2.07%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.9310"
Set-TextValue $ws.Range("E9") "0.55%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1467"
Set-TextValue $ws.Range("E10") "3.52%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1934"
Set-TextValue $ws.Range("E11") "1.27%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.09059"
Set-TextValue $ws.Range("E12") "0.34%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03496"
Set-TextValue $ws.Range("E13") "1.59%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09789"
Set-TextValue $ws.Range("E14") "-1.21%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001390"
Set-TextValue $ws.Range("E15") "-1.00%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.005868"
Set-TextValue $ws.Range("E16") "-2.77%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.763"
Set-TextValue $ws.Range("E17") "-2.02%"
Set-TextValue $ws.Range("E18") "0.72%"
Set-TextValue $ws.Range("D19") "0.3445"
Set-TextValue $ws.Range("E19") "0.00%"
Set-TextValue $ws.Range("D20") "0.1302"
Set-TextValue $ws.Range("E20") "-2.39%"
Set-TextValue $ws.Range("D21") "4.796"
Set-TextValue $ws.Range("E21") "-0.78%"
Set-TextValue $ws.Range("D22") "0.2505"
Set-TextValue $ws.Range("E22") "-4.14%"
Set-TextValue $ws.Range("D23") "0.04376"
Set-TextValue $ws.Range("E23") "0.19%"
Set-TextValue $ws.Range("E24") "0.58%"
Set-TextValue $ws.Range("D25") "0.004283"
Set-TextValue $ws.Range("E25") "-0.29%"
Set-TextValue $ws.Range("D26") "0.0001300"
Set-TextValue $ws.Range("E26") "0.03%"
Set-TextValue $ws.Range("D39") "0.02060"
Set-TextValue $ws.Range("E39") "2.84%"
Set-TextValue $ws.Range("D40") "0.05047"
Set-TextValue $ws.Range("E40") "-1.28%"
Set-TextValue $ws.Range("D41") "0.007437"
Set-TextValue $ws.Range("E41") "-0.89%"
Set-TextValue $ws.Range("D42") "0.01011"
Set-TextValue $ws.Range("E42") "0.06%"
Set-TextValue $ws.Range("D43") "0.1351"
Set-TextValue $ws.Range("E43") "-0.63%"
Set-TextValue $ws.Range("D44") "0.002141"
Set-TextValue $ws.Range("E44") "-1.36%"
Set-TextValue $ws.Range("D45") "0.009074"
Set-TextValue $ws.Range("E45") "-5.79%"
Set-TextValue $ws.Range("D46") "0.00006201"
Set-TextValue $ws.Range("E46") "-0.68%"
Set-TextValue $ws.Range("E47") "0.05%"
Set-TextValue $ws.Range("D48") "0.002798"
Set-TextValue $ws.Range("D49") "0.001600"
Set-TextValue $ws.Range("E49") "27.97%"
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("E50") "0.05%"
Set-TextValue $ws.Range("D51") "0.0002001"
Set-TextValue $ws.Range("E51") "0.05%"
